$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '278.01'
Set-TextValue $ws.Range("E2") '6.49%'
Set-TextValue $ws.Range("D3") '27.45'
Set-TextValue $ws.Range("E3") '2.13%'
Set-TextValue $ws.Range("D4") '4.833'
Set-TextValue $ws.Range("E4") '2.82%'
Set-TextValue $ws.Range("D5") '0.06252'
Set-TextValue $ws.Range("E5") '0.57%'
Set-TextValue $ws.Range("D6") '6.904'
Set-TextValue $ws.Range("E6") '2.22%'
Set-TextValue $ws.Range("B7") 'MXToken'
Set-TextValue $ws.Range("C7") 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range("D7") '0.8777'
Set-TextValue $ws.Range("E7") '3.07%'
Set-TextValue $ws.Range("B8") 'FTXToken'
Set-TextValue $ws.Range("C8") 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range("D8") '0.9437'
Set-TextValue $ws.Range("E8") '3.26%'
Set-TextValue $ws.Range("B9") 'WazirX'
Set-TextValue $ws.Range("C9") 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range("D9") '0.1450'
Set-TextValue $ws.Range("E9") '3.23%'
Set-TextValue $ws.Range("B10") 'LiechtensteinCryptoassetsExchange'
Set-TextValue $ws.Range("C10") 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range("D10") '0.05217'
Set-TextValue $ws.Range("E10") '5.64%'
Set-TextValue $ws.Range("B11") 'MandalaExchangeToken'
Set-TextValue $ws.Range("C11") 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range("D11") '0.07299'
Set-TextValue $ws.Range("E11") '3.27%'
Set-TextValue $ws.Range("B12") 'BitrueCoin'
Set-TextValue $ws.Range("C12") 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range("D12") '0.03160'
Set-TextValue $ws.Range("E12") '1.94%'
Set-TextValue $ws.Range("B13") 'BitMartToken'
Set-TextValue $ws.Range("C13") 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range("D13") '0.09051'
Set-TextValue $ws.Range("E13") '0.04%'
Set-TextValue $ws.Range("B14") 'BitForexToken'
Set-TextValue $ws.Range("C14") 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range("D14") '0.001559'
Set-TextValue $ws.Range("E14") '2.08%'
Set-TextValue $ws.Range("B15") 'One'
Set-TextValue $ws.Range("C15") 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws.Range("D15") '0.0006275'
Set-TextValue $ws.Range("E15") '1.76%'
Set-TextValue $ws.Range("B16") 'TigerCash'
Set-TextValue $ws.Range("C16") 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range("D16") '0.006005'
Set-TextValue $ws.Range("E16") '-0.74%'
Set-TextValue $ws.Range("B17") 'LEO'
Set-TextValue $ws.Range("C17") 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range("D17") '3.460'
Set-TextValue $ws.Range("E17") '0.51%'
Set-TextValue $ws.Range("B18") 'GateToken'
Set-TextValue $ws.Range("C18") 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range("D18") '3.272'
Set-TextValue $ws.Range("E18") '2.90%'
Set-TextValue $ws.Range("D19") '2.286'
Set-TextValue $ws.Range("E19") '6.54%'
Set-TextValue $ws.Range("D21") '0.1312'
Set-TextValue $ws.Range("E21") '0.14%'
Set-TextValue $ws.Range("D22") '3.853'
Set-TextValue $ws.Range("E22") '-6.17%'
Set-TextValue $ws.Range("D23") '0.04329'
Set-TextValue $ws.Range("E23") '2.30%'
Set-TextValue $ws.Range("D24") '0.001177'
Set-TextValue $ws.Range("E24") '-2.48%'
Set-TextValue $ws.Range("D25") '0.004273'
Set-TextValue $ws.Range("E25") '4.86%'
Set-TextValue $ws.Range("D26") '0.0001202'
Set-TextValue $ws.Range("E26") '0.13%'
Set-TextValue $ws.Range("D27") '0.0001692'
Set-TextValue $ws.Range("E27") '3.16%'
Set-TextValue $ws.Range("E40") '1.99%'
Set-TextValue $ws.Range("D41") '0.006104'
Set-TextValue $ws.Range("E41") '47.61%'
Set-TextValue $ws.Range("D42") '0.1153'
Set-TextValue $ws.Range("E42") '3.70%'
Set-TextValue $ws.Range("D43") '0.002134'
Set-TextValue $ws.Range("E43") '-3.47%'
Set-TextValue $ws.Range("D44") '0.01211'
Set-TextValue $ws.Range("E44") '-12.74%'
Set-TextValue $ws.Range("D45") '0.00005075'
Set-TextValue $ws.Range("E45") '-1.68%'
Set-TextValue $ws.Range("D46") '0.00000000751'
Set-TextValue $ws.Range("E46") '0.07%'
Set-TextValue $ws.Range("D47") '2.376'
Set-TextValue $ws.Range("E47") '849.38%'
Set-TextValue $ws.Range("D49") '0.00002102'
Set-TextValue $ws.Range("E49") '0.07%'
Set-TextValue $ws.Range("D50") '0.0002002'
Set-TextValue $ws.Range("E50") '0.07%'
